# Auto-generated Excel COM-interop script to update cryptos.xlsx data
# per commit: "Updated cryptos list on Fri Jun 21 11:59:20 UTC 2024 with GitHub Actions"
#
# Column D (Price) values are stored as text in the original workbook (t="inlineStr"),
# and several of them look like plain numbers (e.g. "1.00", "0.385", "0.0000180"). A
# straight assignment of such a string to Range.Value lets Excel auto-convert it into a
# genuine number cell (dropping trailing zeros, switching to scientific notation, etc.),
# which would not match the source data. To avoid that, values written to column D are
# prefixed with a leading apostrophe - exactly as if a user typed it into the Excel UI -
# which forces text interpretation while leaving the stored text and General number
# format untouched. Columns B, C and E never look numeric, so they are set verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.753.66"
$ws.Range("E2").Value = "  -3.67%  "
$ws.Range("D3").Value = "'3.504.55"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'585.73"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("D6").Value = "'132.07"
$ws.Range("E6").Value = "  -5.20%  "
$ws.Range("D7").Value = "'3.503.73"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "'7.13"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "'0.385"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "'4.095.97"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "'27.84"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.118"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "'3.500.31"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "'63.930.04"
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").Value = "'10.11"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'14.52"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").Value = "'5.70"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "'393.08"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "'0.579"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").Value = "'3.640.13"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "'73.02"
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -6.88%  "
$ws.Range("D28").Value = "'1.60"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.47"
$ws.Range("E29").Value = "  -8.81%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").Value = "'8.26"
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("D33").Value = "'3.507.91"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("D35").Value = "'23.87"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("D36").Value = "'0.144"
$ws.Range("E36").Value = "  -4.25%  "
$ws.Range("D37").Value = "'5.31"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "'1.58"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").Value = "'6.94"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "'167.43"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "'0.0807"
$ws.Range("E41").Value = "  -4.36%  "
$ws.Range("D42").Value = "'26.95"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").Value = "'0.813"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'41.62"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").Value = "'1.21"
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("D47").Value = "'4.40"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("D48").Value = "'1.64"
$ws.Range("E48").Value = "  -4.90%  "
$ws.Range("D49").Value = "'6.92"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").Value = "'2.449.20"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'0.904"
$ws.Range("E51").Value = "  -1.37%  "
